# edit.ps1 - applies the "Did testing (should remember to fix the problems)"
# commit to "Physics todo.docx":
#   1. Adds five new bulleted to-do items above the existing first bullet
#      ("Add contact solver").
#   2. Adds a new "dyn4j" link (and a short comment under it) under the
#      "Links to look at:" heading, just above the existing dyn4j GJK link.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: insert a new ListParagraph-styled bullet immediately before the
# paragraph at $beforeParaIndex. $runs is an array of strings - each element
# becomes (logically) a separate run of text in the new paragraph.
# ---------------------------------------------------------------------------
function Insert-BulletBefore($beforeParaIndex, [string[]]$runs) {
    $beforePara = $d.Paragraphs($beforeParaIndex)
    $r = $beforePara.Range
    $r.Collapse(1)
    $r.InsertParagraphBefore()

    $newPara = $d.Paragraphs($beforeParaIndex)
    $newPara.Range.Text = $runs[0]

    for ($i = 1; $i -lt $runs.Count; $i++) {
        $cur = $d.Paragraphs($beforeParaIndex)
        $endR = $d.Range($cur.Range.Start, $cur.Range.End - 1)
        $endR.Collapse(0)
        $endR.InsertAfter($runs[$i])
    }
}

# Locate the "Add contact solver" bullet (currently the very first bullet).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Add contact solver*") {
        $targetIndex = $i
        break
    }
}

# Inserted in reverse order so each new bullet ends up directly above the
# previous one, producing the right final top-to-bottom order.
Insert-BulletBefore $targetIndex @("Add compound collision properly (also add modifier tool to bind two", " colliders together", ") ")
Insert-BulletBefore $targetIndex @("Make step forward account for iteration count")
Insert-BulletBefore $targetIndex @("Change polygon  pointcase to see if the point is to the left of all lines, if so, return true")
Insert-BulletBefore $targetIndex @("Remove polygon normals, they" + [char]0x2019 + "re not used (alternatively: remember why I calculated normals)")
Insert-BulletBefore $targetIndex @("Make collision points a dequeue", ", take points out of the front when full", " (and also make it inside a collision callback function)")

# ---------------------------------------------------------------------------
# Add the new "dyn4j" link + "^-- broadphase stuff" comment right below the
# "Links to look at:" heading.
# ---------------------------------------------------------------------------
$linksIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "Links to look at:*") {
        $linksIndex = $i
        break
    }
}

$linksPara = $d.Paragraphs($linksIndex)
$lr = $linksPara.Range
$lr.Collapse(0)
$lr.InsertParagraphAfter()

$hyperlinkParaIndex = $linksIndex + 1
$d.Paragraphs($hyperlinkParaIndex).Range.Text = "dyn4j"
$hpPara = $d.Paragraphs($hyperlinkParaIndex)
$textRange = $d.Range($hpPara.Range.Start, $hpPara.Range.End - 1)
$d.Hyperlinks.Add($textRange, "https://dyn4j.org/", $null, $null, "dyn4j") | Out-Null

$hpPara2 = $d.Paragraphs($hyperlinkParaIndex)
$hr = $hpPara2.Range
$hr.Collapse(0)
$hr.InsertParagraphAfter()
$commentParaIndex = $hyperlinkParaIndex + 1
$d.Paragraphs($commentParaIndex).Range.Text = "^-- broadphase stuff"

Write-Output "Done"
